$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A19").Value = "BiB_Metabolomics.metms_pairm_s.metmss42370"
$ws.Range("A20").Value = "BiB_Metabolomics.metms_pairm_s.metmss485"
$ws.Range("A21").Value = "BiB_Metabolomics.metms_pairm_s.*"

$ws.Range("A22").Select()
